# "Data translation & tasks edition"
# Row 15 (task #13 - Patrik Masrna / Filters in Auctions, Bid progress /
# Enabling to select no value or multiple values for filters) is cleared
# out, turning it back into a blank task row like the rows below it.
# Clearing the content also drops the now-unused shared strings
# ("Patrik Masrna", "Filters in Auctions, Bid progress",
# "Enabling to select no value or multiple values for filters") when the
# workbook is saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C15:E15").ClearContents()

# Move the active selection onto D15, matching where editing left off.
$ws.Range("D15").Select()
